$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (years header): add 2021 / 2022 ---
$ws.Range("R4").Value = 2021
$ws.Range("S4").Value = 2022

# --- Row 5 ---
$ws.Range("P5").Value = 25.6
$ws.Range("Q5").Value = 23.8
$ws.Range("R5").Value = 26.8
$ws.Range("S5").Value = 26.8

# --- Row 6 ---
$ws.Range("P6").Value = 18.6
$ws.Range("Q6").Value = 16.7
$ws.Range("R6").Value = 19.3
$ws.Range("S6").Value = 19.3

# --- Row 7 ---
$ws.Range("R7").Value = "-"
$ws.Range("S7").Value = "-"

# --- Row 8 ---
$ws.Range("P8").Value = 2.1
$ws.Range("Q8").Value = 1.8
$ws.Range("R8").Value = 1.8
$ws.Range("S8").Value = 1.8

# --- Row 9 ---
$ws.Range("P9").Value = 4.9
$ws.Range("Q9").Value = 5.2
$ws.Range("R9").Value = 5.7
$ws.Range("S9").Value = 5.7

# --- Row 10 ---
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 0

# Make the new R/S columns pick up the same formatting as column Q
# (same number format / font / borders / alignment as the preceding year column)
# so they render identically to the rest of the table.
$ws.Range("Q4:Q10").Copy()
$ws.Range("R4:R10").PasteSpecial(-4122)

$ws.Range("Q4:Q10").Copy()
$ws.Range("S4:S10").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Selection moves to T3 (matches the recorded selection in the edit) ---
$ws.Range("T3").Select()
